$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = '76.594.50'
$ws.Range("E2").Value = '  +0.97%  '
$ws.Range("D3").Value = '3.052.47'
$ws.Range("E3").Value = '  +5.01%  '
$style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '201.99'
$ws.Range("D5").Style = $style
$ws.Range("E5").Value = '  +0.10%  '
$style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '629.63'
$ws.Range("D6").Style = $style
$ws.Range("E6").Value = '  +5.50%  '
$ws.Range("E8").Value = '  +0.76%  '
$style = $ws.Range("D9").Style
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.210'
$ws.Range("D9").Style = $style
$ws.Range("E9").Value = '  +7.08%  '
$ws.Range("D10").Value = '3.050.49'
$ws.Range("E10").Value = '  +4.97%  '
$style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.441'
$ws.Range("D11").Style = $style
$ws.Range("E11").Value = '  +2.05%  '
$style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.161'
$ws.Range("D12").Style = $style
$ws.Range("E12").Value = '  -0.52%  '
$style = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.16'
$ws.Range("D13").Style = $style
$ws.Range("E13").Value = '  +5.81%  '
$ws.Range("D14").Value = '3.609.70'
$ws.Range("E14").Value = '  +4.71%  '
$style = $ws.Range("D15").Style
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '29.75'
$ws.Range("D15").Style = $style
$ws.Range("E15").Value = '  +7.46%  '
$ws.Range("D16").Value = '76.492.61'
$ws.Range("E16").Value = '  +0.96%  '
$ws.Range("E17").Value = '  +2.97%  '
$ws.Range("D18").Value = '3.045.22'
$ws.Range("E18").Value = '  +4.76%  '
$style = $ws.Range("D19").Style
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '13.53'
$ws.Range("D19").Style = $style
$ws.Range("E19").Value = '  +4.80%  '
$style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '9.06'
$ws.Range("D20").Style = $style
$ws.Range("E20").Value = '  +3.94%  '
$style = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '376.83'
$ws.Range("D21").Style = $style
$ws.Range("E21").Value = '  +1.39%  '
$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.30'
$ws.Range("D22").Style = $style
$ws.Range("E22").Value = '  -0.50%  '
$ws.Range("B23").Value = 'Polkadot'
$ws.Range("C23").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.37'
$ws.Range("D23").Style = $style
$ws.Range("E23").Value = '  +1.64%  '
$ws.Range("E24").Value = '  +3.83%  '
$ws.Range("D25").Value = '3.192.06'
$ws.Range("E25").Value = '  +4.53%  '
$ws.Range("B26").Value = 'NEARProtocol'
$ws.Range("C26").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$style = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '4.41'
$ws.Range("D26").Style = $style
$ws.Range("E26").Value = '  +4.81%  '
$ws.Range("B27").Value = 'Dai'
$ws.Range("C27").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("D27").Style = $style
$ws.Range("E27").Value = '  -0.05%  '
$style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.01'
$ws.Range("D28").Style = $style
$ws.Range("E28").Value = '  +3.55%  '
$ws.Range("E29").Value = '  +4.46%  '
$ws.Range("E30").Value = '  +0.00%  '
$style = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '8.34'
$ws.Range("D31").Style = $style
$ws.Range("E31").Value = '  +8.75%  '
$style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.43'
$ws.Range("D32").Style = $style
$ws.Range("E32").Value = '  +2.33%  '
$style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '515.53'
$ws.Range("D33").Style = $style
$ws.Range("E33").Value = '  +3.17%  '
$ws.Range("E34").Value = '  +8.23%  '
$style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("D35").Style = $style
$ws.Range("E35").Value = '  +0.04%  '
$style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '21.01'
$ws.Range("D36").Style = $style
$ws.Range("E36").Value = '  +4.47%  '
$style = $ws.Range("D37").Style
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.62'
$ws.Range("D37").Style = $style
$ws.Range("E37").Value = '  -0.80%  '
$style = $ws.Range("D38").Style
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.388'
$ws.Range("D38").Style = $style
$ws.Range("E38").Value = '  +11.73%  '
$style = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '20.06'
$ws.Range("D39").Style = $style
$ws.Range("E39").Value = '  +2.27%  '
$ws.Range("E40").Value = '  +3.24%  '
$style = $ws.Range("D41").Style
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '189.36'
$ws.Range("D41").Style = $style
$ws.Range("E41").Value = '  +4.95%  '
$style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.113'
$ws.Range("D42").Style = $style
$ws.Range("E42").Value = '  +0.47%  '
$ws.Range("E43").Value = '  +0.26%  '
$style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '5.23'
$ws.Range("D44").Style = $style
$ws.Range("E44").Value = '  +5.37%  '
$ws.Range("E45").Value = '  +7.11%  '
$style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '42.31'
$ws.Range("D46").Style = $style
$ws.Range("E46").Value = '  +5.77%  '
$ws.Range("E47").Value = '  +2.12%  '
$style = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.737'
$ws.Range("D48").Style = $style
$ws.Range("E48").Value = '  +12.85%  '
$style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.46'
$ws.Range("D49").Style = $style
$ws.Range("E49").Value = '  +5.00%  '
$style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.611'
$ws.Range("D50").Style = $style
$ws.Range("E50").Value = '  +7.07%  '
$ws.Range("E51").Value = '  +5.23%  '
